$wb = $excel.ActiveWorkbook

# Neodymium
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = [double]"0.0002195405251500087"
$ws.Range("C3").Value = [double]"0.01062411525673284"
$ws.Range("C4").Value = [double]"0.009608716352691784"
$ws.Range("C5").Value = [double]"2.138791829054013E-07"

# Dysprosium
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = [double]"0.0002195405251500235"
$ws.Range("C3").Value = [double]"0.01062411525673355"
$ws.Range("C4").Value = [double]"0.009608716352692431"
$ws.Range("C5").Value = [double]"2.138791829054185E-07"

# Copper
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = [double]"0.007632681444695514"
$ws.Range("C3").Value = [double]"0.027535891297259"
$ws.Range("C4").Value = [double]"0.007370778697872926"
$ws.Range("C5").Value = [double]"0.016160587324431"

# Raw silicon
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = [double]"0.005750015024097243"
$ws.Range("C3").Value = [double]"0.01921210602835477"
$ws.Range("C4").Value = [double]"0.0053924808017845"
$ws.Range("C5").Value = [double]"0.006847896595910315"
